# Nov 1st push - Updated 2020 mortality datasets; Pneumonia & Influenza
# added to CCB Cause List; year parameter added; New url parameter tab
#
# Insert five new rows into the "Date | Update" news table, directly
# after the header row, each describing one update item (newest first).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each new row is inserted immediately before the (growing) row index
# right after the header, so the five rows end up in the same
# top-to-bottom order they are added in.

# --- Row 1: 11/01/2021 -- Respiratory Infections disaggregation -------
$row = $t.Rows.Add($t.Rows.Item(2))
$row.Cells.Item(1).Range.Text = "11/01/2021 --"
$row.Cells.Item(2).Range.Text = "“Respiratory Infections” have been disaggregated into “Influenza” and “Pneumonia” (and a small number of other respiratory infections, which are included only in “Other Infectious Diseases/Nutritional Deficiencies”), and congenital pneumonia has been moved into “Neonatal conditions.”"

# --- Row 2: 11/01/2021 -- births checkbox ------------------------------
$row = $t.Rows.Add($t.Rows.Item(3))
$row.Cells.Item(1).Range.Text = "11/01/2021 -- "
$row.Cells.Item(2).Range.Text = "A check box has been added to include or exclude “births” from the hospitalization charts. The default is “off” which excludes births from the charts."

# --- Row 3: 10/01/2021 -- 2020 Patient Discharge data ------------------
$row = $t.Rows.Add($t.Rows.Item(4))
$row.Cells.Item(1).Range.Text = "10/01/2021 --"
$row.Cells.Item(2).Range.Text = "2020 Patient Discharge data added"

# --- Row 4: 10/01/2021 -- State of Public Health Reports ---------------
$row = $t.Rows.Add($t.Rows.Item(5))
$row.Cells.Item(1).Range.Text = "10/01/2021 --"
$row.Cells.Item(2).Range.Text = "2019-2021 State of Public Health Reports posted on the homepage."

# --- Row 5: 10/01/2021 -- Excess Mortality Data Brief updated ----------
$row = $t.Rows.Add($t.Rows.Item(6))
$row.Cells.Item(1).Range.Text = "10/01/2021 --"
$row.Cells.Item(2).Range.Text = "2020 Excess Mortality Data Brief updated with Quarter 1, 2021 data, and new content."

Write-Host "Inserted 5 rows; table now has" $t.Rows.Count "rows"
